$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Four row-pairs had their match-detail/odds columns (F:V) swapped between
#    each other (columns A:E - Indice, pais, torneio, temporada, data_partida
#    - stay on their original row). Do each swap via a scratch range far
#    below the used range so the source values aren't clobbered mid-copy.
# ---------------------------------------------------------------------------
$ws.Range("F58:V58").Copy($ws.Range("F500"))
$ws.Range("F59:V59").Copy($ws.Range("F58"))
$ws.Range("F500:V500").Copy($ws.Range("F59"))
$ws.Range("F500:V500").ClearContents()

$ws.Range("F63:V63").Copy($ws.Range("F501"))
$ws.Range("F65:V65").Copy($ws.Range("F63"))
$ws.Range("F501:V501").Copy($ws.Range("F65"))
$ws.Range("F501:V501").ClearContents()

$ws.Range("F66:V66").Copy($ws.Range("F502"))
$ws.Range("F67:V67").Copy($ws.Range("F66"))
$ws.Range("F502:V502").Copy($ws.Range("F67"))
$ws.Range("F502:V502").ClearContents()

$ws.Range("F88:V88").Copy($ws.Range("F503"))
$ws.Range("F89:V89").Copy($ws.Range("F88"))
$ws.Range("F503:V503").Copy($ws.Range("F89"))
$ws.Range("F503:V503").ClearContents()

# ---------------------------------------------------------------------------
# 2) Append a new match row (row 120) at the end of the table, reusing the
#    formatting of the last existing row (119) then overwriting its values.
# ---------------------------------------------------------------------------
$ws.Range("A119:V119").Copy($ws.Range("A120"))

$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "spain"
$ws.Range("C120").Value = "laliga"
$ws.Range("D120").Value = "2023-2024"
$ws.Range("E120").Value = 45236.875
$ws.Range("F120").Value = "Getafe"
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = "Cadiz CF"
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 1.93
$ws.Range("K120").Value = "22/10/2023 12:02"
$ws.Range("L120").Value = 2.08
$ws.Range("M120").Value = "06/11/2023 20:59"
$ws.Range("N120").Value = 3.19
$ws.Range("O120").Value = "22/10/2023 12:02"
$ws.Range("P120").Value = 3.1
$ws.Range("Q120").Value = "06/11/2023 20:59"
$ws.Range("R120").Value = 4.87
$ws.Range("S120").Value = "22/10/2023 12:02"
$ws.Range("T120").Value = 4.41
$ws.Range("U120").Value = "06/11/2023 20:59"
$ws.Range("V120").Value = "https://www.betexplorer.com/football/spain/laliga/getafe-cadiz/xdbTDHba/"
